$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2662068.2
$ws.Range("C7").Value = -40.08504626517798
$ws.Range("D7").Value = 2680
$ws.Range("E7").Value = 2680
$ws.Range("F7").Value = 993.3090298507464
$ws.Range("G7").Value = 5.879560032879527
